$d = $word.ActiveDocument

# --- Edit 1: "schedule. Blacklisted and removed the Nouveau driver to prevent"
# becomes "schedule. " + bold("Blacklisted and removed the Nouveau driver") + " to prevent"
$full = $d.Content.Text
$target = "Blacklisted and removed the Nouveau driver"
$idx = $full.IndexOf($target)
if ($idx -lt 0) { throw "Could not find target text for bold formatting" }
$r = $d.Range($idx, $idx + $target.Length)
$r.Bold = 1

# --- Edit 2: ", resulting in inaccurate data early on. Configured "
# becomes ", resulting in inaccurate delay data. Configured "
$d.Content.Find.Execute("resulting in inaccurate data early on. Configured", $true, $false, $false, $false, $false,
                         $true, 1, $false, "resulting in inaccurate delay data. Configured", 2)
